$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-02-13 Thursday" "2025-02-14 Friday"

Replace-Text "50×71=" "96×77="
Replace-Text "99×71=" "16×98="
Replace-Text "97×88=" "71×63="
Replace-Text "67×84=" "52×25="
Replace-Text "53×39=" "95×25="
Replace-Text "42×94=" "69×56="
Replace-Text "31×84=" "20×85="
Replace-Text "45×43=" "96×82="
Replace-Text "86×36=" "98×50="
Replace-Text "28×57=" "29×32="
Replace-Text "76×86=" "79×99="
Replace-Text "93×28=" "12×18="
Replace-Text "33×71=" "78×65="
Replace-Text "28×36=" "96×57="
Replace-Text "25×60=" "64×51="
Replace-Text "32×18=" "29×33="
Replace-Text "68×48=" "91×12="
Replace-Text "91×53=" "15×21="
Replace-Text "17×16=" "52×35="
Replace-Text "59×26=" "20×92="
Replace-Text "69×75=" "31×42="
Replace-Text "58×73=" "21×85="
Replace-Text "58×80=" "75×87="
Replace-Text "71×33=" "18×43="
Replace-Text "76×55=" "20×16="
